# Populate the first row of Sheet1 with four header-style text values and
# leave the selection on A2 (as if the user had just typed the last value
# and pressed Enter to drop down to the next row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "data 1"
$ws.Range("B1").Value = "data 2"
$ws.Range("C1").Value = "data 3"
$ws.Range("D1").Value = "data 4"

$ws.Range("A2").Select()
